$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.556.08"
$ws.Range("E2").Value = "  -0.32%  "

# Row 3
$ws.Range("D3").Value = "2.265.61"
$ws.Range("E3").Value = "  -0.79%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "120.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.69%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "264.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.10%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.643"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.97%  "

# Row 8
$ws.Range("E8").Value = "  +0.22%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.620"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.09%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.56"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.23%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0943"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.00%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.07%  "

# Row 13
$ws.Range("E13").Value = "  -1.48%  "

# Row 14
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.912"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.24%  "

# Row 15
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.20%  "

# Row 16
$ws.Range("D16").Value = "2.611.62"
$ws.Range("E16").Value = "  -0.64%  "

# Row 17
$ws.Range("D17").Value = "2.269.57"
$ws.Range("E17").Value = "  -0.79%  "

# Row 18
$ws.Range("D18").Value = "43.533.89"
$ws.Range("E18").Value = "  -0.10%  "

# Row 19
$ws.Range("E19").Value = "  +0.48%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.25%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.27%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.85%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.74%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.81%  "

# Row 25
$ws.Range("E25").Value = "  -0.80%  "

# Row 26
$ws.Range("E26").Value = "  +1.92%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.36%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "41.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.34%  "

# Row 29
$ws.Range("E29").Value = "  -0.90%  "

# Row 30
$ws.Range("E30").Value = "  +0.36%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.02"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.83%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.65"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.42%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0916"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.92%  "

# Row 34
$ws.Range("E34").Value = "  +0.20%  "

# Row 35
$ws.Range("E35").Value = "  +1.65%  "

# Row 36
$ws.Range("E36").Value = "  +14.73%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0375"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.56%  "

# Row 38
$ws.Range("E38").Value = "  -1.94%  "

# Row 39
$ws.Range("E39").Value = "  +0.03%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.54"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.33%  "

# Row 41
$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.39"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.09%  "

# Row 42
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.02%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.237"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.58%  "

# Row 44
$ws.Range("E44").Value = "  -0.01%  "

# Row 45
$ws.Range("E45").Value = "  -1.34%  "

# Row 46
$ws.Range("E46").Value = "  -9.02%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "73.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +40.66%  "

# Row 48
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.52"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.07%  "

# Row 49
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.26"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.15%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.37%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "101.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.39%  "
